$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 contains account 000330949 / RENATO / 3000 - remove it entirely,
# shifting subsequent rows up.
$ws.Rows.Item(3).Delete()
